# Updated cryptos list on Sun Nov 26 10:08:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to keep its text value exactly as given (avoid Excel
    # re-interpreting dotted numbers like "234.11" or "1.00" as a float and
    # mangling/rounding it), without leaving a permanent number-format change
    # on the cell (restore to Normal style afterwards).
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "37.805.40"
Set-TextValue "E2" "  +0.10%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.090.84"
Set-TextValue "E3" "  +0.31%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "234.11"
Set-TextValue "E5" "  -0.09%  "

# Row 6 - XRP
Set-TextValue "E6" "  -0.17%  "

# Row 7 - was USDC, now Solana
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D7" "58.37"
Set-TextValue "E7" "  -0.46%  "

# Row 8 - was Solana, now USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.01%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +0.64%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0782"
Set-TextValue "E10" "  -0.41%  "

# Row 11 - TRON
Set-TextValue "D11" "0.108"
Set-TextValue "E11" "  +2.67%  "

# Row 12 - Chainlink
Set-TextValue "D12" "15.22"
Set-TextValue "E12" "  +1.91%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.398.72"

# Row 14 - Avalanche
Set-TextValue "D14" "21.23"
Set-TextValue "E14" "  +0.74%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.779"
Set-TextValue "E15" "  +0.60%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.37"
Set-TextValue "E16" "  +1.12%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.092.44"
Set-TextValue "E17" "  +0.53%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "37.754.58"
Set-TextValue "E18" "  +0.21%  "

# Row 19 - Uniswap
Set-TextValue "D19" "6.13"
Set-TextValue "E19" "  -0.56%  "

# Row 20 - Litecoin
Set-TextValue "D20" "71.03"
Set-TextValue "E20" "  -0.27%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0838"
Set-TextValue "E21" "  +0.44%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "229.91"
Set-TextValue "E22" "  +0.59%  "

# Row 23 - Dai
Set-TextValue "E23" "  -0.03%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -0.76%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +0.04%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.71"
Set-TextValue "E26" "  +7.97%  "

# Row 27 - Monero
Set-TextValue "D27" "171.58"
Set-TextValue "E27" "  +1.25%  "

# Row 28 - Kaspa
Set-TextValue "E28" "  -2.58%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "19.52"
Set-TextValue "E29" "  -0.02%  "

# Row 30 - ImmutableX
Set-TextValue "E30" "  -0.66%  "

# Row 31 - Stellar
Set-TextValue "E31" "  +0.53%  "

# Row 32 - Filecoin
Set-TextValue "E32" "  +0.32%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0634"
Set-TextValue "E33" "  +0.18%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "4.62"
Set-TextValue "E34" "  -1.06%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "2.51"
Set-TextValue "E35" "  +0.78%  "

# Row 36 - WEMIXToken
Set-TextValue "E36" "  -0.08%  "

# Row 37 - RenderToken
Set-TextValue "E37" "  -2.09%  "

# Row 38 - BinanceUSD
Set-TextValue "E38" "  +0.16%  "

# Row 39 - THORChain
Set-TextValue "D39" "5.38"
Set-TextValue "E39" "  -0.36%  "

# Row 40 - VeChain
Set-TextValue "D40" "0.0238"
Set-TextValue "E40" "  +10.54%  "

# Row 41 - Aave
Set-TextValue "D41" "101.57"
Set-TextValue "E41" "  +3.24%  "

# Row 42 - Cronos
Set-TextValue "D42" "0.0971"
Set-TextValue "E42" "  -0.90%  "

# Row 43 - was HuobiToken, now TrustWalletToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "1.21"
Set-TextValue "E43" "  +4.56%  "

# Row 44 - was TrustWalletToken, now HuobiToken
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D44" "2.91"
Set-TextValue "E44" "  +1.54%  "

# Row 45 - InjectiveProtocol
Set-TextValue "D45" "16.71"
Set-TextValue "E45" "  +0.92%  "

# Row 46 - Maker
Set-TextValue "D46" "1.451.57"
Set-TextValue "E46" "  -0.42%  "

# Row 47 - FTXToken
Set-TextValue "D47" "4.12"
Set-TextValue "E47" "  -4.41%  "

# Row 48 - ARBITRUM
Set-TextValue "E48" "  -0.59%  "

# Row 49 - FraxShare
Set-TextValue "D49" "7.23"
Set-TextValue "E49" "  -2.61%  "

# Row 50 - MXToken
Set-TextValue "E50" "  -1.81%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.281.97"
Set-TextValue "E51" "  +0.26%  "
